$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Candidates")

# --- 1. Swap the "twtrHandle" (E) and "party" (G) columns, row by row,
#        and fill in the new "wikipediaPage" (H) column --------------------

# Header row: E1 <-> G1, new H1
$e1 = $ws.Range("E1").Value()
$g1 = $ws.Range("G1").Value()
$ws.Range("E1").Value = $g1
$ws.Range("G1").Value = $e1
$ws.Range("H1").Value = "wikipediaPage"

# Give H1 the same (header) formatting as the other header cells
$ws.Range("G1").Copy()
$ws.Range("H1").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = 0

$wikiPages = @{
    2  = "Anura_Kumara_Dissanayake"
    3  = "A._S._P._Liyanage"
    4  = "Janaka_Ratnayake"
    5  = "M._A._Sumanthiran"
    6  = "Unknown"
    7  = "Unknown"
    8  = "Ranil_Wickremesinghe"
    9  = "Sarath_Fonseka"
    10 = "Unknown"
    11 = "Sajith_Premadasa"
    12 = "Wijeyadasa_Rajapakshe"
}

for ($r = 2; $r -le 12; $r++) {
    $eVal = $ws.Cells.Item($r, 5).Value()
    $gVal = $ws.Cells.Item($r, 7).Value()
    $ws.Cells.Item($r, 5).Value = $gVal
    $ws.Cells.Item($r, 7).Value = $eVal
    $ws.Cells.Item($r, 8).Value = $wikiPages[$r]
}

# --- 2. Column widths --------------------------------------------------------
# Columns C (firstName), D (lastName) and F (imgFile) keep their original
# (bestFit) widths untouched - only A, B, E, G and the new H need resizing.
$ws.Columns.Item(1).ColumnWidth = 4.05
$ws.Columns.Item(2).ColumnWidth = 2.94
$ws.Columns.Item(5).ColumnWidth = 10.83
$ws.Columns.Item(7).ColumnWidth = 15.39
$ws.Columns.Item(8).ColumnWidth = 25.28

# --- 3. Page setup -----------------------------------------------------------
$ws.PageSetup.Orientation = 1   # xlPortrait

# --- 4. Selection -------------------------------------------------------------
$ws.Activate() | Out-Null
$ws.Range("I21").Select() | Out-Null
